# Refresh cryptos list: updated Price (D) / Volume(1h) (E) values for rows 2-51,
# plus two coin-row swaps (Polkadot/Chainlink at rows 13-14; TrustWalletToken/Hedera at rows 37-38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.655.55"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.843.47"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'312.51"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4265"
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("D8").Value = "'0.3615"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").Value = "'0.07307"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "'0.8697"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").Value = "'20.68"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "1.859.29"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'6.535"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.328"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "'0.06992"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "'79.37"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "'0.000008962"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "'15.28"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").Value = "27.655.42"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'4.976"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'10.34"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "2.108.38"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("D26").Value = "'155.20"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'18.50"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'120.07"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("D29").Value = "'5.221"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "'1.873"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").Value = "'0.08896"
$ws.Range("D32").Value = "'0.7643"
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("D33").Value = "'2.964"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("D34").Value = "'4.496"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").Value = "'1.125"
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05428"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.102"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "'2.815"
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").Value = "'0.1661"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").Value = "'0.5060"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'6.561"
$ws.Range("E43").Value = "  -4.22%  "
$ws.Range("D44").Value = "'8.399"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").Value = "'0.06547"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "'106.16"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").Value = "'10.31"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "'0.4635"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("D50").Value = "'1.632"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "'64.40"
$ws.Range("E51").Value = "  +0.08%  "
